# This script updates the "F" column (view/attendance counts) values for a
# handful of expo rows on the "展览" (Exhibition) sheet and on the
# "全部类型" (All Types) summary sheet, which mirrors the same rows.
#
# Mapping of old -> new values (row indices differ between the two sheets
# because "全部类型" also interleaves rows from the "演出" sheet):
#   5294  -> 5305
#   10803 -> 10842
#   269   -> 270
#   576   -> 577
#   151   -> 154
#   196   -> 206
#   914   -> 920

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5305
$wsExpo.Range("F4").Value = 10842
$wsExpo.Range("F5").Value = 270
$wsExpo.Range("F6").Value = 577
$wsExpo.Range("F7").Value = 154
$wsExpo.Range("F8").Value = 206
$wsExpo.Range("F9").Value = 920

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5305
$wsAll.Range("F7").Value = 10842
$wsAll.Range("F8").Value = 270
$wsAll.Range("F9").Value = 577
$wsAll.Range("F10").Value = 154
$wsAll.Range("F13").Value = 206
$wsAll.Range("F14").Value = 920
